$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: fill in "Implemented" column G with "yes"
$ws.Range("G12").Value = "yes"

# Row 13: previously an empty placeholder row, now populated with new test case data
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "CoreDataTests.ComposedFoodItemBehavior"
$ws.Range("C13").Value = "ComposedFoodItem.create(from composedFoodItemVM: ComposedFoodItemViewModel, _ isImport: Bool)"
$ws.Range("E13").Value = "Empty DB"
$ws.Range("D13").Value = "isImport = true"
$ws.Range("F13").Value = "ComposedFoodItem, its related FoodItem, all related Ingredients and their related FoodItems are created"
$ws.Range("G13").Value = "yes"

# Adjust row height for row 13 to match wrapped multi-line content
$ws.Rows.Item(13).RowHeight = 34

# Column width adjustments (col B widened to fit the new longer "Test suite" text,
# col F widened slightly for the new "Expected result" text)
$ws.Columns.Item(2).ColumnWidth = 37.3
$ws.Columns.Item(6).ColumnWidth = 43.8

# Update the active selection to reflect the new cursor position
$ws.Range("G14").Select()
